# edit.ps1
# Applies the commit: "updated templates, added instructions, changed some layout formatting"
#
# Per the OOXML diff this corresponds to:
#  1. Adding a new worksheet named "Assay" (sheetId 5) as the last tab,
#     containing a two-column Setting/Value table of assay parameters.
#  2. Moving the active/selected tab from "Compounds" to "Patterns", and
#     updating the Patterns sheet's selected cell to K21.
#     (Compounds keeps its own last-selected cell, B37, it just stops being
#     the tab that is active/open when the workbook is reopened.)

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Assay" worksheet at the end, populate with settings ---

$sheets = $wb.Worksheets
$afterSheet = $sheets.Item($sheets.Count)

$assay = $sheets.Add($null, $afterSheet)
$assay.Name = "Assay"

$assayData = @(
    @("Setting", "Value"),
    @("DMSO Tolerance", 0.005),
    @("Well Volume (µL)", 25),
    @("Backfill (µL)", 10),
    @("Allowed Error", 0.1),
    @("Destination Replicates", 1),
    @("Use Intermediate Plates", 1),
    @("DMSO Normalization", 1)
)

for ($i = 0; $i -lt $assayData.Length; $i++) {
    $rowNum = $i + 1
    $assay.Cells.Item($rowNum, 1).Value = $assayData[$i][0]
    $assay.Cells.Item($rowNum, 2).Value = $assayData[$i][1]
}

$assay.Range("A1:B8").Select() | Out-Null

# --- 2. Make "Patterns" the active tab with K21 selected ---

$patterns = $wb.Worksheets.Item("Patterns")
$patterns.Activate()
$patterns.Range("K21").Select() | Out-Null
